# Auto-generated edit script: applies the value changes described by the
# Bahamut_Profits.xlsx diff (per-sheet Leve-profit recalculation update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 323.58
$ws.Range("I15").Value = 323.58
$ws.Range("K15").Value = 970.74
$ws.Range("M15").Value = -801.74
$ws.Range("H40").Value = 2293.6428
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2316.2307
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2316.2307
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2666.2307
$ws.Range("H74").Value = 4632598.5
$ws.Range("I74").Value = 5095358
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 5095358
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -5094422
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 4632598.5
$ws.Range("I77").Value = 5095358
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 25476790
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -25472110
$ws.Range("N77").Value = -34360
$ws.Range("H103").Value = 1281.8182
$ws.Range("I103").Value = 1683.3334
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 5050.0002
$ws.Range("L103").Value = 2400
$ws.Range("M103").Value = -4464.0002
$ws.Range("N103").Value = -3572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1062.2106
$ws.Range("I61").Value = 860.1667
$ws.Range("J61").Value = 1408.5714
$ws.Range("K61").Value = 860.1667
$ws.Range("L61").Value = 1408.5714
$ws.Range("M61").Value = -648.1667
$ws.Range("N61").Value = -1832.5714
$ws.Range("H88").Value = 3409.6365
$ws.Range("I88").Value = 3063.25
$ws.Range("J88").Value = 4333.3335
$ws.Range("K88").Value = 3063.25
$ws.Range("L88").Value = 4333.3335
$ws.Range("M88").Value = -2657.25
$ws.Range("N88").Value = -5145.3335
$ws.Range("H91").Value = 3409.6365
$ws.Range("I91").Value = 3063.25
$ws.Range("J91").Value = 4333.3335
$ws.Range("K91").Value = 3063.25
$ws.Range("L91").Value = 4333.3335
$ws.Range("M91").Value = -1659.25
$ws.Range("N91").Value = -7141.3335
$ws.Range("H102").Value = 4504.5557
$ws.Range("I102").Value = 4716.25
$ws.Range("J102").Value = 2811
$ws.Range("K102").Value = 4716.25
$ws.Range("L102").Value = 2811
$ws.Range("M102").Value = -3094.25
$ws.Range("N102").Value = -6055
$ws.Range("H132").Value = 1523.4
$ws.Range("I132").Value = 993.86957
$ws.Range("J132").Value = 2538.3333
$ws.Range("K132").Value = 2981.60871
$ws.Range("L132").Value = 7614.999899999999
$ws.Range("M132").Value = -451.60871
$ws.Range("N132").Value = -12674.9999
$ws.Range("H136").Value = 1062.2106
$ws.Range("I136").Value = 860.1667
$ws.Range("J136").Value = 1408.5714
$ws.Range("K136").Value = 2580.5001
$ws.Range("L136").Value = 4225.7142
$ws.Range("M136").Value = -30.5001000000002
$ws.Range("N136").Value = -9325.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 282.5
$ws.Range("I22").Value = 280
$ws.Range("J22").Value = 285
$ws.Range("K22").Value = 280
$ws.Range("L22").Value = 285
$ws.Range("M22").Value = -107
$ws.Range("N22").Value = -631
$ws.Range("H35").Value = 26900
$ws.Range("I35").Value = 20000
$ws.Range("J35").Value = 33800
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 33800
$ws.Range("M35").Value = -19690
$ws.Range("N35").Value = -34420
$ws.Range("H86").Value = 2243
$ws.Range("I86").Value = 2097.6365
$ws.Range("J86").Value = 2699.8572
$ws.Range("K86").Value = 2097.6365
$ws.Range("L86").Value = 2699.8572
$ws.Range("M86").Value = -974.6365000000001
$ws.Range("N86").Value = -4945.8572
$ws.Range("H89").Value = 2243
$ws.Range("I89").Value = 2097.6365
$ws.Range("J89").Value = 2699.8572
$ws.Range("K89").Value = 10488.1825
$ws.Range("L89").Value = 13499.286
$ws.Range("M89").Value = -4872.182500000001
$ws.Range("N89").Value = -24731.286
$ws.Range("H99").Value = 2477.7778
$ws.Range("I99").Value = 2475
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2475
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -977
$ws.Range("N99").Value = -5496
$ws.Range("H134").Value = 57369.25
$ws.Range("I134").Value = 2290.6843
$ws.Range("J134").Value = 118927.65
$ws.Range("K134").Value = 6872.0529
$ws.Range("L134").Value = 356782.95
$ws.Range("M134").Value = -4337.0529
$ws.Range("N134").Value = -361852.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5850.826
$ws.Range("I58").Value = 2012.7273
$ws.Range("J58").Value = 9369.083000000001
$ws.Range("K58").Value = 2012.7273
$ws.Range("L58").Value = 9369.083000000001
$ws.Range("M58").Value = -1809.7273
$ws.Range("N58").Value = -9775.083000000001
$ws.Range("H62").Value = 9800
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 9800
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H136").Value = 5850.826
$ws.Range("I136").Value = 2012.7273
$ws.Range("J136").Value = 9369.083000000001
$ws.Range("K136").Value = 6038.1819
$ws.Range("L136").Value = 28107.249
$ws.Range("M136").Value = -3488.1819
$ws.Range("N136").Value = -33207.249

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1737.5
$ws.Range("I81").Value = 975
$ws.Range("K81").Value = 2925
$ws.Range("M81").Value = -1802
$ws.Range("H84").Value = 1737.5
$ws.Range("I84").Value = 975
$ws.Range("K84").Value = 8775
$ws.Range("M84").Value = -3159
$ws.Range("H141").Value = 7007.722
$ws.Range("I141").Value = 5353.5454
$ws.Range("J141").Value = 9607.143
$ws.Range("K141").Value = 16060.6362
$ws.Range("L141").Value = 28821.429
$ws.Range("M141").Value = -10880.6362
$ws.Range("N141").Value = -39181.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 33348860
$ws.Range("J42").Value = 33348860
$ws.Range("L42").Value = 33348860
$ws.Range("N42").Value = -33349830
$ws.Range("H70").Value = 4870.9
$ws.Range("I70").Value = 4100
$ws.Range("J70").Value = 6669.6665
$ws.Range("K70").Value = 4100
$ws.Range("L70").Value = 6669.6665
$ws.Range("M70").Value = -3830
$ws.Range("N70").Value = -7209.6665
$ws.Range("H73").Value = 4870.9
$ws.Range("I73").Value = 4100
$ws.Range("J73").Value = 6669.6665
$ws.Range("K73").Value = 4100
$ws.Range("L73").Value = 6669.6665
$ws.Range("M73").Value = -3164
$ws.Range("N73").Value = -8541.666499999999
$ws.Range("H115").Value = 33348860
$ws.Range("J115").Value = 33348860
$ws.Range("L115").Value = 33348860
$ws.Range("N115").Value = -33351210
$ws.Range("H122").Value = 550163.2
$ws.Range("I122").Value = 823911.4399999999
$ws.Range("J122").Value = 2666.75
$ws.Range("K122").Value = 2471734.32
$ws.Range("L122").Value = 8000.25
$ws.Range("M122").Value = -2469284.32
$ws.Range("N122").Value = -12900.25
$ws.Range("H132").Value = 2997.9443
$ws.Range("I132").Value = 2900.6365
$ws.Range("J132").Value = 3150.8572
$ws.Range("K132").Value = 8701.9095
$ws.Range("L132").Value = 9452.571599999999
$ws.Range("M132").Value = -6171.9095
$ws.Range("N132").Value = -14512.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1387
$ws.Range("I16").Value = 1523.1
$ws.Range("K16").Value = 1523.1
$ws.Range("M16").Value = -1353.1
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H40").Value = 3368171.2
$ws.Range("I40").Value = 3368171.2
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3368171.2
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3368035.2
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 1142.4375
$ws.Range("I46").Value = 1780
$ws.Range("J46").Value = 759.9
$ws.Range("K46").Value = 1780
$ws.Range("L46").Value = 759.9
$ws.Range("M46").Value = -1592
$ws.Range("N46").Value = -1135.9
$ws.Range("H55").Value = 490
$ws.Range("I55").Value = 490
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 490
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -317
$ws.Range("N55").ClearContents()
$ws.Range("H100").Value = 1869.8
$ws.Range("I100").Value = 1924.75
$ws.Range("J100").Value = 1650
$ws.Range("K100").Value = 1924.75
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -1383.75
$ws.Range("N100").Value = -2732

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 70007
$ws.Range("J15").Value = 70007
$ws.Range("L15").Value = 70007
$ws.Range("N15").Value = -70583

